# Update handback report timestamps ("Generate Report for Handback")
$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (D5) and Correspond Handback DateTime (G5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-23 09:22:56"
$wsZhCn.Range("G5").Value = "2016-02-23 09:23:46"

# de-de sheet: Correspond Handoff Datetime (D5) and Correspond Handback DateTime (G5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-23 09:23:08"
$wsDeDe.Range("G5").Value = "2016-02-23 09:24:08"
